# Auto-generated Excel COM-interop script to apply the Shiva_Profits market-data refresh
# (author commit: "chore: update Sheets via scheduled runner")
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2643.0833
$ws.Range("I40").Value = 2449.5
$ws.Range("J40").Value = 3611
$ws.Range("K40").Value = 2449.5
$ws.Range("L40").Value = 3611
$ws.Range("M40").Value = -2274.5
$ws.Range("N40").Value = -3961
$ws.Range("H76").Value = 5100.1113
$ws.Range("I76").Value = 4128.857
$ws.Range("J76").Value = 8499.5
$ws.Range("K76").Value = 4128.857
$ws.Range("L76").Value = 8499.5
$ws.Range("M76").Value = -3813.857
$ws.Range("N76").Value = -9129.5
$ws.Range("H79").Value = 5100.1113
$ws.Range("I79").Value = 4128.857
$ws.Range("J79").Value = 8499.5
$ws.Range("K79").Value = 4128.857
$ws.Range("L79").Value = 8499.5
$ws.Range("M79").Value = -3036.857
$ws.Range("N79").Value = -10683.5
$ws.Range("H92").Value = 15873356
$ws.Range("I92").Value = 23809762
$ws.Range("J92").Value = 543.4286
$ws.Range("K92").Value = 23809762
$ws.Range("L92").Value = 543.4286
$ws.Range("M92").Value = -23808514
$ws.Range("N92").Value = -3039.4286
$ws.Range("H99").Value = 3553.8333
$ws.Range("I99").Value = 1580.75
$ws.Range("J99").Value = 7500
$ws.Range("K99").Value = 4742.25
$ws.Range("L99").Value = 22500
$ws.Range("M99").Value = -3244.25
$ws.Range("N99").Value = -25496
$ws.Range("H101").Value = 1976.7778
$ws.Range("I101").Value = 1968.2
$ws.Range("K101").Value = 5904.6
$ws.Range("M101").Value = -4282.6
$ws.Range("H132").Value = 7532.25
$ws.Range("I132").Value = 4596.7393
$ws.Range("J132").Value = 30037.834
$ws.Range("K132").Value = 13790.2179
$ws.Range("L132").Value = 90113.50199999999
$ws.Range("M132").Value = -11260.2179
$ws.Range("N132").Value = -95173.50199999999
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 180000
$ws.Range("J136").Value = 180000
$ws.Range("L136").Value = 180000
$ws.Range("N136").Value = -190200
$ws.Range("H137").Value = 6705.0884
$ws.Range("I137").Value = 7625.5
$ws.Range("J137").Value = 3713.75
$ws.Range("K137").Value = 22876.5
$ws.Range("L137").Value = 11141.25
$ws.Range("M137").Value = -20326.5
$ws.Range("N137").Value = -16241.25
$ws.Range("H138").Value = 24392674
$ws.Range("J138").Value = 4794.643
$ws.Range("L138").Value = 14383.929
$ws.Range("N138").Value = -24663.929
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4482.974
$ws.Range("I32").Value = 4542.533
$ws.Range("K32").Value = 4542.533
$ws.Range("M32").Value = -4255.533
$ws.Range("H61").Value = 3495.149
$ws.Range("I61").Value = 3440
$ws.Range("J61").Value = 4304
$ws.Range("K61").Value = 3440
$ws.Range("L61").Value = 4304
$ws.Range("M61").Value = -3228
$ws.Range("N61").Value = -4728
$ws.Range("H74").Value = 1827.5952
$ws.Range("I74").Value = 1628.2927
$ws.Range("J74").Value = 9999
$ws.Range("K74").Value = 1628.2927
$ws.Range("L74").Value = 9999
$ws.Range("M74").Value = -754.2927
$ws.Range("N74").Value = -11747
$ws.Range("H77").Value = 1827.5952
$ws.Range("I77").Value = 1628.2927
$ws.Range("J77").Value = 9999
$ws.Range("K77").Value = 8141.4635
$ws.Range("L77").Value = 49995
$ws.Range("M77").Value = -3773.4635
$ws.Range("N77").Value = -58731
$ws.Range("H132").Value = 1488.9524
$ws.Range("I132").Value = 1330.5143
$ws.Range("K132").Value = 3991.5429
$ws.Range("M132").Value = -1461.5429
$ws.Range("H136").Value = 3495.149
$ws.Range("I136").Value = 3440
$ws.Range("J136").Value = 4304
$ws.Range("K136").Value = 10320
$ws.Range("L136").Value = 12912
$ws.Range("M136").Value = -7770
$ws.Range("N136").Value = -18012

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2638.8276
$ws.Range("I31").Value = 1951.35
$ws.Range("J31").Value = 4166.5557
$ws.Range("K31").Value = 1951.35
$ws.Range("L31").Value = 4166.5557
$ws.Range("M31").Value = -1656.35
$ws.Range("N31").Value = -4756.5557
$ws.Range("H34").Value = 2638.8276
$ws.Range("I34").Value = 1951.35
$ws.Range("J34").Value = 4166.5557
$ws.Range("K34").Value = 1951.35
$ws.Range("L34").Value = 4166.5557
$ws.Range("M34").Value = -1749.35
$ws.Range("N34").Value = -4570.5557
$ws.Range("H132").Value = 7054.424
$ws.Range("I132").Value = 3671.6223
$ws.Range("J132").Value = 17927.715
$ws.Range("K132").Value = 11014.8669
$ws.Range("L132").Value = 53783.145
$ws.Range("M132").Value = -8484.866900000001
$ws.Range("N132").Value = -58843.145
$ws.Range("H134").Value = 3130.093
$ws.Range("I134").Value = 2962.7896
$ws.Range("J134").Value = 4401.6
$ws.Range("K134").Value = 8888.3688
$ws.Range("L134").Value = 13204.8
$ws.Range("M134").Value = -6353.3688
$ws.Range("N134").Value = -18274.8

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1222.4
$ws.Range("I97").Value = 765.2
$ws.Range("J97").Value = 1679.6
$ws.Range("K97").Value = 2295.6
$ws.Range("L97").Value = 5038.799999999999
$ws.Range("M97").Value = -1799.6
$ws.Range("N97").Value = -6030.799999999999
$ws.Range("H134").Value = 1715.96
$ws.Range("I134").Value = 1537.4584
$ws.Range("K134").Value = 4612.3752
$ws.Range("M134").Value = 457.6247999999996

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2577.5625
$ws.Range("I132").Value = 1978.6522
$ws.Range("J132").Value = 4108.1113
$ws.Range("K132").Value = 5935.9566
$ws.Range("L132").Value = 12324.3339
$ws.Range("M132").Value = -3405.9566
$ws.Range("N132").Value = -17384.3339
$ws.Range("H141").Value = 53929
$ws.Range("J141").Value = 53929
$ws.Range("L141").Value = 53929
$ws.Range("N141").Value = -64289

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1770.6666
$ws.Range("I61").Value = 1736
$ws.Range("J61").Value = 1996
$ws.Range("K61").Value = 1736
$ws.Range("L61").Value = 1996
$ws.Range("M61").Value = -1534
$ws.Range("N61").Value = -2400
$ws.Range("H68").Value = 20843302
$ws.Range("I68").Value = 41670356
$ws.Range("J68").Value = 16250.25
$ws.Range("K68").Value = 41670356
$ws.Range("L68").Value = 16250.25
$ws.Range("M68").Value = -41669607
$ws.Range("N68").Value = -17748.25
$ws.Range("H71").Value = 20843302
$ws.Range("I71").Value = 41670356
$ws.Range("J71").Value = 16250.25
$ws.Range("K71").Value = 208351780
$ws.Range("L71").Value = 81251.25
$ws.Range("M71").Value = -208348036
$ws.Range("N71").Value = -88739.25
$ws.Range("H113").Value = 1770.6666
$ws.Range("I113").Value = 1736
$ws.Range("J113").Value = 1996
$ws.Range("K113").Value = 1736
$ws.Range("L113").Value = 1996
$ws.Range("M113").Value = 434
$ws.Range("N113").Value = -6336

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3194.054
$ws.Range("I132").Value = 2822.3462
$ws.Range("K132").Value = 8467.0386
$ws.Range("M132").Value = -5937.0386
$ws.Range("H136").Value = 2686.3438
$ws.Range("I136").Value = 2482.6775
$ws.Range("J136").Value = 9000
$ws.Range("K136").Value = 7448.032499999999
$ws.Range("L136").Value = 27000
$ws.Range("M136").Value = -4898.032499999999
$ws.Range("N136").Value = -32100
